$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.984.63"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").Value = "3.220.24"
$ws.Range("E3").Value = "  -3.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.46"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.08"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "3.226.80"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.609"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.72"
$ws.Range("E11").Value = "  -8.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.135"
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.11"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "3.734.78"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.116"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.215.81"
$ws.Range("E17").Value = "  -4.12%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.35"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").Value = "62.803.31"
$ws.Range("E19").Value = "  -2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.11"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.970"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "367.39"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.18"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.21"
$ws.Range("E25").Value = "  +4.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("E27").Value = "  +2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.67"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.37"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.58"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "644.05"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.53"
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.35"
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("E35").Value = "  +2.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.04"
$ws.Range("E36").Value = "  -4.34%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.02"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").Value = "0.0₃0726"
$ws.Range("E40").Value = "  +17.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "2.889.29"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.52"
$ws.Range("E44").Value = "  +9.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.95"
$ws.Range("E45").Value = "  +11.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.68"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0394"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.01"
$ws.Range("E49").Value = "  +8.05%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.124"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.86"
$ws.Range("E51").Value = "  -0.19%  "
